# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-09 13:19:59
#
# A new ANATOMY attendance session (session 2) was recorded for 19 Year-2/C1
# students: 19 rows are appended to the end of the "Attendance" sheet, and
# the corresponding "Summary" sheet metrics (risk %, sessions needed,
# attended/missed totals, attended-ANATOMY totals, and in a few cases the
# risk Status label) are refreshed to reflect the new attendance record.

$wb = $excel.ActiveWorkbook

$newRows = @(
    @{ R=333; ID="202107"; Name="خديجة ادم محمد ادريس موسى"; Time="14:37:29" },
    @{ R=334; ID="210926"; Name="كوثر محمد المبارك يوسف"; Time="14:37:32" },
    @{ R=335; ID="212412"; Name="شارلز ماكوك مالوك"; Time="14:37:37" },
    @{ R=336; ID="210870"; Name="عبد الله جعفر عثمان جعفر"; Time="14:37:44" },
    @{ R=337; ID="212442"; Name="رميساء محى الدين الامين الطيب"; Time="14:37:48" },
    @{ R=338; ID="200540"; Name="محمد سعيد ابراهيم عواد درويش"; Time="14:37:52" },
    @{ R=339; ID="220370"; Name="اغيد مرزوق الرفاعى"; Time="14:37:56" },
    @{ R=340; ID="211704"; Name="الياس احمد بكردان"; Time="14:38:11" },
    @{ R=341; ID="212211"; Name="سوزان جيريمياه لادو"; Time="14:38:19" },
    @{ R=342; ID="211294"; Name="نور الهدى خلاوى الشحاذه"; Time="14:38:26" },
    @{ R=343; ID="221319"; Name="روان صلاح طاهر الوهباني"; Time="14:38:47" },
    @{ R=344; ID="220811"; Name="يعقوب يوسف يوسف"; Time="14:38:53" },
    @{ R=345; ID="220766"; Name="ميار بنت خالد بن محمد الشيخ"; Time="14:39:01" },
    @{ R=346; ID="212322"; Name="مهند حافظ عابدين الفاضل"; Time="14:39:08" },
    @{ R=347; ID="220314"; Name="احمد ربيع قطب عبد المطلب بهوت"; Time="14:39:17" },
    @{ R=348; ID="220428"; Name="بسمله محمد عبد الحميد محمد"; Time="14:39:27" },
    @{ R=349; ID="221682"; Name="سرين حاج صدوق"; Time="14:39:34" },
    @{ R=350; ID="201441"; Name="هند محمد ادم عيسى"; Time="14:39:41" },
    @{ R=351; ID="210998"; Name="زينب نادر عوض السيد عبد القادر"; Time="14:39:56" }
)

$summaryChanges = @(
    @{ R=3; G="17.2%"; H=17; L=5; M=0; O=2; F=$null },
    @{ R=5; G="6.9%"; H=20; L=2; M=3; O=2; F="Low Risk" },
    @{ R=12; G="10.3%"; H=19; L=3; M=2; O=2; F=$null },
    @{ R=19; G="10.3%"; H=19; L=3; M=2; O=2; F=$null },
    @{ R=20; G="10.3%"; H=19; L=3; M=2; O=1; F=$null },
    @{ R=24; G="13.8%"; H=18; L=4; M=1; O=1; F="No Risk" },
    @{ R=32; G="6.9%"; H=20; L=2; M=3; O=1; F="Low Risk" },
    @{ R=38; G="10.3%"; H=19; L=3; M=2; O=1; F=$null },
    @{ R=58; G="3.4%"; H=21; L=1; M=4; O=1; F=$null },
    @{ R=64; G="17.2%"; H=17; L=5; M=0; O=2; F=$null },
    @{ R=68; G="10.3%"; H=19; L=3; M=2; O=2; F=$null },
    @{ R=71; G="10.3%"; H=19; L=3; M=2; O=1; F=$null },
    @{ R=76; G="10.3%"; H=19; L=3; M=2; O=2; F=$null },
    @{ R=77; G="10.3%"; H=19; L=3; M=2; O=1; F=$null },
    @{ R=78; G="10.3%"; H=19; L=3; M=2; O=1; F=$null },
    @{ R=84; G="13.8%"; H=18; L=4; M=1; O=1; F="No Risk" },
    @{ R=85; G="10.3%"; H=19; L=3; M=2; O=2; F=$null },
    @{ R=110; G="6.9%"; H=20; L=2; M=3; O=2; F="Low Risk" },
    @{ R=176; G="17.2%"; H=17; L=5; M=0; O=2; F=$null }
)

# --- 1) Append the 19 new Attendance rows (333-351) ---
$attendance = $wb.Worksheets.Item(2)

foreach ($row in $newRows) {
    $r = $row.R

    $attendance.Range("A$r").NumberFormat = "@"
    $attendance.Range("A$r").Value = $row.ID
    $attendance.Range("A$r").Style = "Normal"

    $attendance.Range("B$r").Value = $row.Name
    $attendance.Range("C$r").Value = "Year 2"
    $attendance.Range("D$r").Value = "C1"
    $attendance.Range("E$r").Value = "$($row.ID)@med.asu.edu.eg"
    $attendance.Range("F$r").Value = "ANATOMY"

    $attendance.Range("G$r").NumberFormat = "@"
    $attendance.Range("G$r").Value = "2"
    $attendance.Range("G$r").Style = "Normal"

    $attendance.Range("H$r").Value = "ANATOMY"

    $attendance.Range("I$r").NumberFormat = "@"
    $attendance.Range("I$r").Value = "09/11/2025"
    $attendance.Range("I$r").Style = "Normal"

    $attendance.Range("J$r").Value = $row.Time
    $attendance.Range("K$r").Value = "C1"
}

# Re-sync the AutoFilter range + the hidden _xlnm._FilterDatabase defined
# name to cover the newly appended rows (1..351).
$attendance.Range("A1:K351").AutoFilter() | Out-Null
$attendance.Range("A1:K351").AutoFilter() | Out-Null

for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "Attendance!_FilterDatabase") {
        $n.RefersTo = "='Attendance'!`$A`$1:`$K`$351"
    }
}

# --- 2) Refresh the affected Summary-sheet metrics ---
$summary = $wb.Worksheets.Item(1)

# Reference cells already carrying the target "Status" styles, used as
# formatting donors for rows whose risk level changes (copy format only).
$lowRiskDonor = $summary.Range("F2")    # style for "Low Risk"
$noRiskDonor  = $summary.Range("F22")   # style for "No Risk"

foreach ($chg in $summaryChanges) {
    $r = $chg.R

    if ($chg.F) {
        if ($chg.F -eq "Low Risk") {
            $lowRiskDonor.Copy() | Out-Null
        } elseif ($chg.F -eq "No Risk") {
            $noRiskDonor.Copy() | Out-Null
        }
        $summary.Range("F$r").PasteSpecial(-4122) | Out-Null
        $summary.Range("F$r").Value = $chg.F
    }

    $summary.Range("G$r").NumberFormat = "@"
    $summary.Range("G$r").Value = $chg.G
    $summary.Range("G$r").NumberFormat = "0.0%"

    $summary.Range("H$r").Value = $chg.H
    $summary.Range("L$r").Value = $chg.L
    $summary.Range("M$r").Value = $chg.M
    $summary.Range("O$r").Value = $chg.O
}
$excel.CutCopyMode = 0
